$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "2025-06-19" / BEMOL S/A rows (id_venda 365782 and 366707).
# This shifts every following row up by two (old row 4 -> new row 2, etc.)
$ws.Range("A2:A3").EntireRow.Delete()

# Corrected figures on the rows that shifted up (values per the updated dataset)
$ws.Range("G2").Value = -89    # id_venda 370495: estoque_atualizado -88 -> -89
$ws.Range("G3").Value = 9      # id_venda 374463: estoque_atualizado -455 -> 9
$ws.Range("G5").Value = -142   # id_venda 374491: estoque_atualizado -140 -> -142
$ws.Range("H5").Value = 1.03   # id_venda 374491: media_vendas 1.04 -> 1.03
$ws.Range("I5").Value = 0.18   # id_venda 374491: desvio_padrao 0.19 -> 0.18
